# Update relay settings rows 2-20 with the new calculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "R" column (TCC time-dial style value) per row, taken from the updated
# calculation - these replace the old R values (which previously all sat at 0.7).
$rValues = @{
    2  = 1.9712896934131849
    3  = 1.7039851694832704
    4  = 0.80187537387448016
    5  = 0.99752525475971998
    6  = 0.50117210867155004
    7  = 0.10023442173431002
    8  = 0.75175816300732523
    9  = 0.66501683650647991
    10 = 0.50117210867155004
    11 = 0.20046884346862004
    12 = 1.0023442173431001
    13 = 0.79042001139055906
    14 = 1.4128571428571426
    15 = 0.40093768693724008
    16 = 0.40093768693724008
    17 = 1.4128571428571426
    18 = 1.4128571428571426
    19 = 1.3832350199334786
    20 = 1.976050028476398
}

foreach ($row in 2..20) {
    $ws.Range("F$row").Value = 13
    $ws.Range("N$row").Value = 0.7
    $ws.Range("Q$row").Value = 0
    $ws.Range("R$row").Value = $rValues[$row]
    $ws.Range("R$row").NumberFormat = "0.00"
    $ws.Range("S$row").Value = 1
    $ws.Range("T$row").Value = 1
    $ws.Range("U$row").Value = 0.5
    $ws.Range("V$row").Value = 1.2
}

# Update the active selection on the sheet to match the author's final cursor position.
$ws.Range("G24").Select()
